$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:60 down to 44:61
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43 with the new data record
$ws.Range("A43").Value = 9
$ws.Range("B43").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C43").Value = "Metropolitana"
$ws.Range("D43").Value = 44455
$ws.Range("D43").NumberFormat = $ws.Range("D44").NumberFormat
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 100112022
$ws.Range("G43").Value = "Arveja Verde"
$ws.Range("H43").Value = "Perfection"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 18
$ws.Range("K43").Value = 32000
$ws.Range("L43").Value = 33000
$ws.Range("M43").Value = 32500
$ws.Range("N43").Value = "`$/malla 25 kilos"
$ws.Range("O43").Value = "Provincia de Huasco"
$ws.Range("P43").Value = 1300
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"
